$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.799.48"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "3.065.25"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'537.28"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("D6").Value = "'133.10"
$ws.Range("E6").Value = "  -3.78%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.058.40"
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").Value = "'0.493"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("D11").Value = "'6.09"
$ws.Range("E11").Value = "  -9.05%  "
$ws.Range("D12").Value = "'0.452"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "'0.0000224"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "'34.25"
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("D15").Value = "3.550.60"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "62.739.21"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "3.059.98"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "'6.64"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("D20").Value = "'482.07"
$ws.Range("E20").Value = "  -4.32%  "
$ws.Range("D21").Value = "'13.30"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").Value = "'0.695"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").Value = "'7.10"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").Value = "'79.08"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "'12.09"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("D28").Value = "'8.07"
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "'25.95"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'1.87"
$ws.Range("E31").Value = "  -8.39%  "
$ws.Range("D32").Value = "'1.11"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("E33").Value = "  -6.86%  "
$ws.Range("D34").Value = "'56.99"
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("D35").Value = "'5.33"
$ws.Range("E35").Value = "  +2.87%  "
$ws.Range("D36").Value = "'6.00"
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("D37").Value = "'484.34"
$ws.Range("E37").Value = "  -9.30%  "
$ws.Range("D38").Value = "3.112.26"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("D40").Value = "'0.0795"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("D42").Value = "'8.07"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "'2.61"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "'0.252"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "0.0₃0533"
$ws.Range("E46").Value = "  +6.69%  "
$ws.Range("D47").Value = "'121.70"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'2.01"
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D49").Value = "'24.40"
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "'2.29"
$ws.Range("E51").Value = "  -2.36%  "
